$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.61900842859776
$ws.Range("C2").Value = 10.44134461977274
$ws.Range("D2").Value = 9.642584506989472
$ws.Range("F2").Value = 33.38749113950155
$ws.Range("G2").Value = 34.71411329326438
$ws.Range("H2").Value = 15.910238371016
$ws.Range("J2").Value = 10.76767496243022
$ws.Range("K2").Value = 9.53196645767798
$ws.Range("L2").Value = 11.31174274511446
$ws.Range("N2").Value = 19.80595528425099
$ws.Range("O2").Value = 24.98518049869629
$ws.Range("B3").Value = 13.40365289755233
$ws.Range("C3").Value = 10.45588753714687
$ws.Range("D3").Value = 9.625505030781763
$ws.Range("F3").Value = 33.46196501129079
$ws.Range("G3").Value = 34.80760788843919
$ws.Range("H3").Value = 15.95449018266681
$ws.Range("J3").Value = 10.79013303547321
$ws.Range("K3").Value = 9.370433840715968
$ws.Range("L3").Value = 11.30629027134798
$ws.Range("N3").Value = 19.86024246409501
$ws.Range("O3").Value = 25.05983583396824
$ws.Range("B4").Value = 13.27183586264525
$ws.Range("C4").Value = 10.46546748509545
$ws.Range("D4").Value = 9.616441042074747
$ws.Range("F4").Value = 33.51403042430401
$ws.Range("G4").Value = 34.87323800233859
$ws.Range("H4").Value = 15.98370787206132
$ws.Range("J4").Value = 10.80485898692041
$ws.Range("K4").Value = 9.27104140863168
$ws.Range("L4").Value = 11.3043612128114
$ws.Range("N4").Value = 19.89519104575642
$ws.Range("O4").Value = 25.10985561328345
$ws.Range("B5").Value = 13.21828947072973
$ws.Range("C5").Value = 10.4695354059024
$ws.Range("D5").Value = 9.613108226767526
$ws.Range("F5").Value = 33.53684004442364
$ws.Range("G5").Value = 34.90204631529464
$ws.Range("H5").Value = 15.99612946769572
$ws.Range("J5").Value = 10.81109592469584
$ws.Range("K5").Value = 9.2305347069661
$ws.Range("L5").Value = 11.30393353993677
$ws.Range("N5").Value = 19.90984035894559
$ws.Range("O5").Value = 25.13129004562271
$ws.Range("B6").Value = 13.20941035873905
$ws.Range("C6").Value = 10.47022079979975
$ws.Range("D6").Value = 9.612576692969411
$ws.Range("F6").Value = 33.5407237094431
$ws.Range("G6").Value = 34.90695440566724
$ws.Range("H6").Value = 15.99822319098722
$ws.Range("J6").Value = 10.81214583089208
$ws.Range("K6").Value = 9.223809851123772
$ws.Range("L6").Value = 11.30388421997417
$ws.Range("N6").Value = 19.91229751010281
$ws.Range("O6").Value = 25.13491268016189
$ws.Range("B7").Value = 13.27111294032673
$ws.Range("C7").Value = 10.46552168183313
$ws.Range("D7").Value = 9.616394629713392
$ws.Range("F7").Value = 33.51433159647762
$ws.Range("G7").Value = 34.8736181724785
$ws.Range("H7").Value = 15.98387330763452
$ws.Range("J7").Value = 10.80494214428114
$ws.Range("K7").Value = 9.270495069215103
$ws.Range("L7").Value = 11.3043539916652
$ws.Range("N7").Value = 19.89538696015044
$ws.Range("O7").Value = 25.11014043076843
$ws.Range("B8").Value = 13.54470802974622
$ws.Range("C8").Value = 10.44622425748154
$ws.Range("D8").Value = 9.636402076465723
$ws.Range("F8").Value = 33.41185356601713
$ws.Range("G8").Value = 34.74464073550278
$ws.Range("H8").Value = 15.92507182782674
$ws.Range("J8").Value = 10.77522440045247
$ws.Range("K8").Value = 9.476343209377076
$ws.Range("L8").Value = 11.30956930141112
$ws.Range("N8").Value = 19.82433877127643
$ws.Range("O8").Value = 25.01005363692833
$ws.Range("B9").Value = 14.08161146193712
$ws.Range("C9").Value = 10.41352461315637
$ws.Range("D9").Value = 9.686783772806503
$ws.Range("F9").Value = 33.26123310405132
$ws.Range("G9").Value = 34.55715683797312
$ws.Range("H9").Value = 15.82598374163278
$ws.Range("J9").Value = 10.72435812810725
$ws.Range("K9").Value = 9.876210301209579
$ws.Range("L9").Value = 11.33097639061952
$ws.Range("N9").Value = 19.6977842542851
$ws.Range("O9").Value = 24.84697625211592
$ws.Range("B10").Value = 14.47252696565956
$ws.Range("C10").Value = 10.39260907606816
$ws.Range("D10").Value = 9.730390210226881
$ws.Range("F10").Value = 33.18131504274394
$ws.Range("G10").Value = 34.45953887628563
$ws.Range("H10").Value = 15.76304552372537
$ws.Range("J10").Value = 10.69147462581174
$ws.Range("K10").Value = 10.16491463612066
$ws.Range("L10").Value = 11.35340961577562
$ws.Range("N10").Value = 19.61251658011832
$ws.Range("O10").Value = 24.74742078180166
$ws.Range("B11").Value = 14.64882681157927
$ws.Range("C11").Value = 10.38376360199958
$ws.Range("D11").Value = 9.751610819478081
$ws.Range("F11").Value = 33.15164045006444
$ws.Range("G11").Value = 34.4238836966877
$ws.Range("H11").Value = 15.73654899463493
$ws.Range("J11").Value = 10.67748336711874
$ws.Range("K11").Value = 10.29459819933063
$ws.Range("L11").Value = 11.36504495066839
$ws.Range("N11").Value = 19.57538504925426
$ws.Range("O11").Value = 24.70653336400434
$ws.Range("B12").Value = 14.71530597176252
$ws.Range("C12").Value = 10.38050983517034
$ws.Range("D12").Value = 9.759841047937991
$ws.Range("F12").Value = 33.14136415531912
$ws.Range("G12").Value = 34.41164301004942
$ws.Range("H12").Value = 15.72682195674008
$ws.Range("J12").Value = 10.6723239090316
$ws.Range("K12").Value = 10.34342526145389
$ws.Range("L12").Value = 11.36965416631504
$ws.Range("N12").Value = 19.56156145424963
$ws.Range("O12").Value = 24.69168353068453
$ws.Range("B13").Value = 14.70100198292955
$ws.Range("C13").Value = 10.38120633661643
$ws.Range("D13").Value = 9.758059947665224
$ws.Range("F13").Value = 33.14353460648627
$ws.Range("G13").Value = 34.41422313374319
$ws.Range("H13").Value = 15.72890321600028
$ws.Range("J13").Value = 10.67342892722772
$ws.Range("K13").Value = 10.33292266932194
$ws.Range("L13").Value = 11.36865249485073
$ws.Range("N13").Value = 19.56452807395869
$ws.Range("O13").Value = 24.69485353227437
$ws.Range("B14").Value = 14.65430208380877
$ws.Range("C14").Value = 10.38349399447644
$ws.Range("D14").Value = 9.752284055168348
$ws.Range("F14").Value = 33.15077575670762
$ws.Range("G14").Value = 34.42285135915265
$ws.Range("H14").Value = 15.7357426019023
$ws.Range("J14").Value = 10.67705611743697
$ws.Range("K14").Value = 10.29862111340324
$ws.Range("L14").Value = 11.36542009673665
$ws.Range("N14").Value = 19.57424302334815
$ws.Range("O14").Value = 24.70529896338453
$ws.Range("B15").Value = 14.62565856318603
$ws.Range("C15").Value = 10.38490771817913
$ws.Range("D15").Value = 9.748771338061216
$ws.Range("F15").Value = 33.15533629880847
$ws.Range("G15").Value = 34.42830070641424
$ws.Range("H15").Value = 15.73997184573916
$ws.Range("J15").Value = 10.67929593043962
$ws.Range("K15").Value = 10.27757251622674
$ws.Range("L15").Value = 11.36346654091075
$ws.Range("N15").Value = 19.58022458598942
$ws.Range("O15").Value = 24.71177958589343
$ws.Range("B16").Value = 14.46096923395666
$ws.Range("C16").Value = 10.39320057897446
$ws.Range("D16").Value = 9.72903085564954
$ws.Range("F16").Value = 33.18338879850734
$ws.Range("G16").Value = 34.46204537315703
$ws.Range("H16").Value = 15.7648200578265
$ws.Range("J16").Value = 10.69240842286163
$ws.Range("K16").Value = 10.15640249992438
$ws.Range("L16").Value = 11.35267778961317
$ws.Range("N16").Value = 19.61497647459736
$ws.Range("O16").Value = 24.75018147162118
$ws.Range("B17").Value = 14.35950110907327
$ws.Range("C17").Value = 10.39845907356088
$ws.Range("D17").Value = 9.71727197129047
$ws.Range("F17").Value = 33.2023093635913
$ws.Range("G17").Value = 34.48499022226082
$ws.Range("H17").Value = 15.78061008401835
$ws.Range("J17").Value = 10.70070004539392
$ws.Range("K17").Value = 10.08161409474492
$ws.Range("L17").Value = 11.34642384011089
$ws.Range("N17").Value = 19.63671936489001
$ws.Range("O17").Value = 24.7748673258741
$ws.Range("B18").Value = 14.30099889738637
$ws.Range("C18").Value = 10.40154662369271
$ws.Range("D18").Value = 9.710639109830762
$ws.Range("F18").Value = 33.21382078376649
$ws.Range("G18").Value = 34.49901111869534
$ws.Range("H18").Value = 15.78989300526816
$ws.Range("J18").Value = 10.70556026801934
$ws.Range("K18").Value = 10.038445029882
$ws.Range("L18").Value = 11.34296152411229
$ws.Range("N18").Value = 19.64938133853914
$ws.Range("O18").Value = 24.7894801791588
$ws.Range("B19").Value = 14.2811688604011
$ws.Range("C19").Value = 10.40260284803135
$ws.Range("D19").Value = 9.708415891973848
$ws.Range("F19").Value = 33.21782634158833
$ws.Range("G19").Value = 34.50389972208903
$ws.Range("H19").Value = 15.79307055954072
$ws.Range("J19").Value = 10.70722151517261
$ws.Range("K19").Value = 10.02380386961143
$ws.Range("L19").Value = 11.34181246819152
$ws.Range("N19").Value = 19.65369529526967
$ws.Range("O19").Value = 24.79449896379314
$ws.Range("B20").Value = 14.37031755417212
$ws.Range("C20").Value = 10.39789278011861
$ws.Range("D20").Value = 9.718510246706002
$ws.Range("F20").Value = 33.20023015527463
$ws.Range("G20").Value = 34.48246244070513
$ws.Range("H20").Value = 15.77890841691269
$ws.Range("J20").Value = 10.69980796245334
$ws.Range("K20").Value = 10.08959158368989
$ws.Range("L20").Value = 11.34707565121125
$ws.Range("N20").Value = 19.63438865491925
$ws.Range("O20").Value = 24.77219660229762
$ws.Range("B21").Value = 14.66802710281899
$ws.Range("C21").Value = 10.3828194564558
$ws.Range("D21").Value = 9.753975334927659
$ws.Range("F21").Value = 33.1486227787993
$ws.Range("G21").Value = 34.42028279410331
$ws.Range("H21").Value = 15.73372538919382
$ws.Range("J21").Value = 10.67598696198822
$ws.Range("K21").Value = 10.30870428658529
$ws.Range("L21").Value = 11.36636403724605
$ws.Range("N21").Value = 19.57138307451305
$ws.Range("O21").Value = 24.70221369448371
$ws.Range("B22").Value = 14.86092838331329
$ws.Range("C22").Value = 10.37352651190732
$ws.Range("D22").Value = 9.778284934689907
$ws.Range("F22").Value = 33.1204946935454
$ws.Range("G22").Value = 34.38699664912538
$ws.Range("H22").Value = 15.70598276294812
$ws.Range("J22").Value = 10.6612269844996
$ws.Range("K22").Value = 10.45024864086581
$ws.Range("L22").Value = 11.38015306248598
$ws.Range("N22").Value = 19.53158805880459
$ws.Range("O22").Value = 24.66016763577375
$ws.Range("B23").Value = 14.75814579681692
$ws.Range("C23").Value = 10.37843537260386
$ws.Range("D23").Value = 9.765208502925869
$ws.Range("F23").Value = 33.13499477178092
$ws.Range("G23").Value = 34.40408865044542
$ws.Range("H23").Value = 15.72062610356375
$ws.Range("J23").Value = 10.66903082709264
$ws.Range("K23").Value = 10.37486950540855
$ws.Range("L23").Value = 11.37268622608408
$ws.Range("N23").Value = 19.55270120720623
$ws.Range("O23").Value = 24.68227045226114
$ws.Range("B24").Value = 14.3654279542607
$ws.Range("C24").Value = 10.39814860084742
$ws.Range("D24").Value = 9.717950024979114
$ws.Range("F24").Value = 33.20116819134583
$ws.Range("G24").Value = 34.48360266683223
$ws.Range("H24").Value = 15.77967710194952
$ws.Range("J24").Value = 10.70021098263096
$ws.Range("K24").Value = 10.08598549396738
$ws.Range("L24").Value = 11.3467805524419
$ws.Range("N24").Value = 19.6354418651452
$ws.Range("O24").Value = 24.77340272624431
$ws.Range("B25").Value = 13.93672030543415
$ws.Range("C25").Value = 10.42182290029675
$ws.Range("D25").Value = 9.671981001998342
$ws.Range("F25").Value = 33.29658411621249
$ws.Range("G25").Value = 34.60084460505751
$ws.Range("H25").Value = 15.85105590651946
$ws.Range("J25").Value = 10.73732855743109
$ws.Range("K25").Value = 9.768733539632956
$ws.Range("L25").Value = 11.32399914561757
$ws.Range("N25").Value = 19.73066104350883
$ws.Range("O25").Value = 24.84697625211592
